# Add "Episode 16" sheet to the AYTO VIP S2 results workbook, as a copy
# of the "Episode 15" sheet (same contestants / match-rate data), placed
# after the last sheet.

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("Episode 15")
$last = $wb.Worksheets.Item($wb.Worksheets.Count)

$source.Copy([System.Reflection.Missing]::Value, $last)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Episode 16"

# Copying a sheet makes it the active one; restore the original active tab
# so the only substantive change is the appended sheet.
$wb.Worksheets.Item(1).Activate()
